# The commit trims a pile of redundant explicit defaults out of
# <w:docDefaults> in word/styles.xml:
#   - rPrDefault/rPr keeps only rFonts, sz, szCs, lang
#     (b/i/smallCaps/strike/color/u/shd/vertAlign are dropped)
#   - pPrDefault/pPr keeps only a bare <w:spacing w:line="276"
#     w:lineRule="auto"/> (keepNext/keepLines/widowControl/pBdr/shd/
#     ind/contextualSpacing/jc are dropped, and spacing's before/after
#     attributes are dropped too)
#
# docDefaults has no Word object-model surface (it isn't a Style object
# Word exposes via the Styles collection), so this is done by round-
# tripping the package text through Document.WordOpenXML: find the
# <w:docDefaults>...</w:docDefaults> block with a regex (robust to the
# exact attribute-order/self-closing quirks WordOpenXML's serializer
# uses) and strip the unwanted child elements in place.

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

$blockMatch = [regex]::Match($xml, '(?s)<w:docDefaults>.*?</w:docDefaults>')
if (-not $blockMatch.Success) {
    throw "docDefaults block not found in WordOpenXML; cannot apply edit"
}
$block = $blockMatch.Value
$newBlock = $block

# --- rPrDefault/rPr: drop these self-closing elements entirely ---
$rprDrop = @('b', 'i', 'smallCaps', 'strike', 'color', 'u', 'shd', 'vertAlign')
foreach ($tag in $rprDrop) {
    $newBlock = [regex]::Replace($newBlock, '<w:' + $tag + '(?:\s[^>]*)?/>', '')
}

# --- pPrDefault/pPr: drop pBdr (and its children) ---
$newBlock = [regex]::Replace($newBlock, '(?s)<w:pBdr>.*?</w:pBdr>', '')

# --- pPrDefault/pPr: drop these other self-closing elements entirely ---
$pprDrop = @('keepNext', 'keepLines', 'widowControl', 'ind', 'contextualSpacing', 'jc')
foreach ($tag in $pprDrop) {
    $newBlock = [regex]::Replace($newBlock, '<w:' + $tag + '(?:\s[^>]*)?/>', '')
}

# --- pPrDefault/pPr: simplify <w:spacing .../> down to line/lineRule only ---
$newBlock = [regex]::Replace($newBlock, '<w:spacing[^/]*/>', '<w:spacing w:line="276" w:lineRule="auto"/>')

$xml = $xml.Substring(0, $blockMatch.Index) + $newBlock + $xml.Substring($blockMatch.Index + $blockMatch.Length)

$d.WordOpenXML = $xml
